# "final edits to playgrounds" -- wording/grammar touch-ups to the Playground
# description document.
#
# The document ships with TrackRevisions turned on (w:trackRevisions in
# settings.xml). Leaving that on would make every Find/Replace below wrap its
# output in w:ins/w:del revision markup instead of plainly rewriting the run
# text, so switch it off first.
$d = $word.ActiveDocument
$d.TrackRevisions = $false

$wdFindContinue = 1
$wdReplaceOne = 1

# --- 1. "I created the ... (joystick and button with a label) programmatically"
#        -> "I programmatically created the ... (specifically the joystick and
#        the button with a label)". i.e. move "programmatically" earlier in the
#        sentence and call out which UI elements were made that way.
$found1 = $d.Content.Find.Execute(
    "I created the game scene and the user interface (joystick and button with a label) programmatically so it’s less intensive on the device running the project",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "I programmatically created the game scene and the user interface (specifically the joystick and the button with a label) so it’s less intensive on the device running the project",
    $wdReplaceOne)

# --- 2. Fix the "bee worried" typo -> "be worried", lower-case the "You" that
#        now sits mid-sentence inside the parenthetical, and tack on the new
#        closing sentences about future work / a fully functioning game.
$found2 = $d.Content.Find.Execute(
    "been caught (You can see this by the caught player being transparent). You will be able to free that player and then try out the different hiding spots on the map. The seeker and the other bots will be frozen, so you don’t have to bee worried about being caught. ",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "been caught (you can see this by the caught player being transparent). You will be able to free that player and then try out the different hiding spots on the map. The seeker and the other bots will be frozen, so you don’t have to be worried about being caught. However, In the future I would want to look into how I could add bot movement to the game, and I would also implement a countdown timer to trigger a victory/loss scene. By doing this, I would have created a fully functioning game.",
    $wdReplaceOne)

# --- 3. Add the missing comma: "student community and I became" -> "student
#        community, and I became".
$found3 = $d.Content.Find.Execute(
    "the student community and I became",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "the student community, and I became",
    $wdReplaceOne)

Write-Output "replace1=$found1 replace2=$found2 replace3=$found3"
